$d = $word.ActiveDocument

function Insert-ParaXml($AfterParagraph, $InnerXml) {
    $AfterParagraph.Range.InsertParagraphAfter()
    $idx = $AfterParagraph.Index + 1
    $newPara = $d.Paragraphs.Item($idx)
    $wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $InnerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($wrapper)
    return $d.Paragraphs.Item($idx)
}

# Locate the "Sprint 2" paragraph (last item of that sprint's bullet list, ilvl 0).
$sprint2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Sprint 2") {
        $sprint2 = $cand
    }
}

# 1) Deadline: Sunday 10th  of May  (ilvl 1)
$p1Xml = '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Deadline: Sunday </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>10</w:t></w:r>' +
    '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">  of</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> May</w:t></w:r>' +
    '</w:p>'
$lastPara = Insert-ParaXml $sprint2 $p1Xml

# 2) Parsa: 1A & 1B and creation/validation on server (ilvl 2)
$p2Xml = '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Parsa</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>: 1A &amp; 1B and creation/validation on server</w:t></w:r>' +
    '</w:p>'
$lastPara = Insert-ParaXml $lastPara $p2Xml

# 3) Jonas: ViewPager Design, ListView -> get info from local storage (ilvl 2)
$p3Xml = '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Jonas: </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ViewPager</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Design, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ListView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> -&gt; get info from local storage</w:t></w:r>' +
    '</w:p>'
$lastPara = Insert-ParaXml $lastPara $p3Xml

# 4) Kristian: Last Views and UX (ilvl 2)
$p4Xml = '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Kristian: Last Views and UX</w:t></w:r>' +
    '</w:p>'
$lastPara = Insert-ParaXml $lastPara $p4Xml

# 5) Morten: Finish Web service and database (ilvl 2) + the relocated _GoBack bookmark
$p5Xml = '<w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Morten: Finish Web service and database</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
$lastPara = Insert-ParaXml $lastPara $p5Xml

# Remove the stale _GoBack bookmark from the old empty ListParagraph further down
# (the one that used to immediately follow "Sprint 4" -- now shifted after our
# 5 freshly-inserted paragraphs). Locate "Sprint 4" then take the paragraph
# right after it.
$sprint4 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Sprint 4") {
        $sprint4 = $cand
    }
}
$emptyPara = $sprint4.Next()
$emptyXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'
$wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $emptyXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$emptyPara.Range.InsertXML($wrapper)

Write-Output "Done"
